$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultipleChoice")
$lo = $ws.ListObjects.Item("MultipleChoiceData")

# New question rows (QuestionID, QuestionText, A, B, C, D, Answer)
$rows = @(
    @(4,  "1+1=", 1, 2, 3, 4,  "A"),
    @(5,  "2+2=", 1, 2, 3, 4,  "D"),
    @(6,  "3+3=", 3, 6, 9, 12, "C"),
    @(7,  "4+4=", 4, 8, 12, 16, "A"),
    @(8,  "5+5=", 5, 10, 15, 25, "D"),
    @(9,  "6+6=", 6, 12, 15, 36, "C"),
    @(10, "7+7=", 7, 3, 7, 49,  "A")
)

foreach ($rowData in $rows) {
    $newRow = $lo.ListRows.Add()
    $dst = $newRow.Range

    # Best-effort format carry-over from the last existing data row
    $ws.Range("A4:G4").Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null

    $dst.Item(1,1).Value = $rowData[0]
    $dst.Item(1,2).Value = $rowData[1]
    $dst.Item(1,3).Value = $rowData[2]
    $dst.Item(1,4).Value = $rowData[3]
    $dst.Item(1,5).Value = $rowData[4]
    $dst.Item(1,6).Value = $rowData[5]
    $dst.Item(1,7).Value = $rowData[6]
}

$ws.Range("C14").Select()
